$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# --- Header row additions (J1 "RGB", K1 "Hex"), matching the style of the
#     existing header row cells (A1:I1) ---
$ws.Range("J1").Value = "RGB"
$ws.Range("K1").Value = "Hex"
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 3: the Good/Moderate "R" color value was corrected from 0 to 255 ---
$ws.Range("G3").Value = 255

# --- Column J: concatenate the R,G,B columns into a "r,g,b" string.
#     Rows 3-6 share one formula (entered once on J3 and filled down to J6). ---
$ws.Range("J2").Formula = "=CONCAT(G2,"","",H2,"","",I2)"
$ws.Range("J3:J6").Formula = "=CONCAT(G3,"","",H3,"","",I3)"
$ws.Range("J7").Formula = "=CONCAT(G7,"","",H7,"","",I7)"

# --- Column K: the matching hex color codes for each AQI band ---
$ws.Range("K2").Value = "#009966"
$ws.Range("K3").Value = "#ffdd33"
$ws.Range("K4").Value = "#ff9933"
$ws.Range("K5").Value = "#cc0033"
$ws.Range("K6").Value = "#660099"
$ws.Range("K7").Value = "#7e0035"

# --- Final view state: scrolled so column E is left-most, K7 selected ---
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("K7").Select()
